# "Fix the dataset name" - update the SAMPLIB/load-library bullet text in the
# install-diagram "Installation Output" box (shape "Rectangle 39" on slide 1).
#
#   Before: "A PDS SAMPLIB member SZWESAMPE containing example JCL"
#   After : "A SAMPLIB library SZWESAMP containing sample members"
#
#   Before: "A PDS load library SZWEAUTH containing load modules"
#   After : "A load library SZWEAUTH containing load modules"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 39")
$tr = $shp.TextFrame.TextRange

# ---------------------------------------------------------------------------
# Paragraph 3: "A PDS SAMPLIB member SZWESAMPE containing example JCL"
# ---------------------------------------------------------------------------
$par3 = $tr.Paragraphs(3, 1)
$base3 = $par3.Start

# Work right-to-left so earlier character offsets (computed against the
# original text) stay valid while later ones are being rewritten.

# "JCL" -> "members"
$run = $tr.Characters($base3 + 50, 3)
$run.Text = "members"

# " " stays " " (kept as its own run, matching the target run split)
$run = $tr.Characters($base3 + 49, 1)
$run.Text = " "

# "example" -> "sample"
$run = $tr.Characters($base3 + 42, 7)
$run.Text = "sample"

# "containing " stays "containing " (kept as its own run)
$run = $tr.Characters($base3 + 31, 11)
$run.Text = "containing "

# "SZWESAMPE " -> "SZWESAMP " (drop the trailing E - the actual dataset name fix)
$run = $tr.Characters($base3 + 21, 10)
$run.Text = "SZWESAMP "

# "A PDS SAMPLIB member " -> "A SAMPLIB library " (rewritten, then re-split below)
$run = $tr.Characters($base3 + 0, 21)
$run.Text = "A SAMPLIB library "

# Re-split the rewritten leading text into its target run boundaries:
# "A SAMPLIB " + "library" + " "
$run = $tr.Characters($base3 + 0, 10)
$run.Text = "A SAMPLIB "

$run = $tr.Characters($base3 + 10, 7)
$run.Text = "library"

$run = $tr.Characters($base3 + 17, 1)
$run.Text = " "

# ---------------------------------------------------------------------------
# Paragraph 4: "A PDS load library SZWEAUTH containing load modules"
# ---------------------------------------------------------------------------
$par4 = $tr.Paragraphs(4, 1)
$base4 = $par4.Start

# "A PDS load library " -> "A load library "
$run = $tr.Characters($base4 + 0, 19)
$run.Text = "A load library "

Write-Host "Paragraph 3:" $tr.Paragraphs(3, 1).Text
Write-Host "Paragraph 4:" $tr.Paragraphs(4, 1).Text
